# "Create and done testimonial section"
# Fill in Start/Finish dates and completion % for the Testinomial section
# (rows 39-41: header "Testinomial" + sub-items "Item" / "Carousel").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 39: Testinomial (section header) ---------------------------------
# Copy the date-formatted style already used on another section header
# (B4 uses cellXf 24: bold font, white fill, thin border, m/d/yyyy format)
# so the new date value picks up matching formatting instead of a fresh style.
$ws.Range("B4").Copy()
$ws.Range("B39").PasteSpecial(-4122)
$ws.Range("B39").Value = 44476
$ws.Range("D39").Value = 0.67

# --- Row 40: Item (sub-item, fully done) -----------------------------------
# Copy style from B5/C5 (cellXf 25: regular font, white fill, thin border,
# m/d/yyyy format) used by other completed sub-item rows.
$ws.Range("B5").Copy()
$ws.Range("B40").PasteSpecial(-4122)
$ws.Range("B40").Value = 44476

$ws.Range("C5").Copy()
$ws.Range("C40").PasteSpecial(-4122)
$ws.Range("C40").Value = 44476

$ws.Range("D40").Value = 1

# --- Row 41: Carousel (sub-item, fully done) -------------------------------
$ws.Range("B5").Copy()
$ws.Range("B41").PasteSpecial(-4122)
$ws.Range("B41").Value = 44476

$ws.Range("C5").Copy()
$ws.Range("C41").PasteSpecial(-4122)
$ws.Range("C41").Value = 44476

$ws.Range("D41").Value = 1

# --- Update the sheet's last selection, matching where editing finished ---
$ws.Range("F40").Select() | Out-Null
$excel.CutCopyMode = $false
